$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Top 3 / Top 4 product identities and update sales figures.
# Row 4 = "Top 3", Row 5 = "Top 4"

# Row 4 (Top 3) now holds the item that used to be "Top 4" (HH0013S / Thắt lưng da)
$ws.Range("C4").Value = "HH0013S"
$ws.Range("D4").Value = "Thắt lưng da"
$ws.Range("F4").Value = 265000
$ws.Range("G4").Value = 342
$ws.Range("H4").Value = 13

# Row 5 (Top 4) now holds the item that used to be "Top 3" (HH0011M / Mũ lưỡi trai)
$ws.Range("C5").Value = "HH0011M"
$ws.Range("D5").Value = "Mũ lưỡi trai"
$ws.Range("F5").Value = 52000
$ws.Range("G5").Value = 110
$ws.Range("H5").Value = 321
